$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestingFile3")

$formula = "=(TestingFile3[[#This Row],[Obj. LR]]-TestingFile3[[#This Row],[LB Heuristic]])/TestingFile3[[#This Row],[Obj. LR]]"

for ($r = 2; $r -le 121; $r++) {
    $ws.Cells.Item($r, 7).Formula = $formula
}

$ws.Range("G121").Select()
